# Weekly refresh of the "Perejil" (parsley) price series for Vega Modelo
# de Temuco: a new sample (dated 2021-08-06 / serial 44414) is inserted at
# the top of the data block (row 39), pushing every subsequent row (40..156)
# down by one and appending the former last row as a new row 157.
#
# Columns A, B, C, E, F, G, H, I, R are identical on every data row, so the
# "shift" only has to touch D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de
# comercializacion), O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) -
# but copying the full row is simplest and harmless since the constant
# columns copy onto themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 39
$lastDataRow  = 156
$newLastRow   = $lastDataRow + 1   # 157
$lastCol      = 18                 # column R

# 1) Append a new last row (157) that is a copy of the current last row (156)
#    before anything else gets overwritten.
for ($col = 1; $col -le $lastCol; $col++) {
    $v = $ws.Cells.Item($lastDataRow, $col).Value2
    $ws.Cells.Item($newLastRow, $col).Value = $v
}
# Preserve the date number format on the newly created row's Fecha cell.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastDataRow, 4).NumberFormat

# 2) Shift rows 156 down through 40 down by one: row r becomes the old
#    row (r - 1). Walk from the bottom up so we never clobber a row before
#    it has been read.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $v = $ws.Cells.Item($r - 1, $col).Value2
        $ws.Cells.Item($r, $col).Value = $v
    }
}

# 3) Write the brand-new sample into row 39 (Volumen/Unidad/Kg-o-Unidades
#    are unchanged from what used to be there).
$ws.Cells.Item($firstDataRow, 4).Value  = 44414                        # Fecha
$ws.Cells.Item($firstDataRow, 11).Value = 3500                         # Precio minimo
$ws.Cells.Item($firstDataRow, 12).Value = 4000                         # Precio maximo
$ws.Cells.Item($firstDataRow, 13).Value = 3750                         # Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 15).Value = "Región de La Araucanía"     # Origen
$ws.Cells.Item($firstDataRow, 16).Value = 1250                         # Precio $/Kg
